# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the 6d83c270-... row on both the zh-cn and de-de sheets,
# simulating a later CI run that re-generated the handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-23 14:48:16"
$wsZhCn.Range("H3").Value = "2016-03-23 14:48:48"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-23 14:48:20"
$wsDeDe.Range("H3").Value = "2016-03-23 14:48:55"
